# map_layers_inventory.xlsx update
# - adds a new "USDM Weeks in Drought" layer entry in row 31 (Sheet1)
# - adds a hyperlink for the new Info Links cell (H31)
# - updates the sheet view's active cell / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data row (row 31) -------------------------------------------------
# Columns: A Layer Name | B Description | C Hosting Agency | D Geoserver Type
#          E Geoserver Link | F Layer ID | G Notes | H Info Links
$ws.Range("B31").Value = "USDM Weeks in Drought"
$ws.Range("C31").Value = "NDMC"
$ws.Range("D31").Value = "csv download"
$ws.Range("E31").Value = "new REST Service available"
$ws.Range("F31").Value = "n/a"
$ws.Range("G31").Value = "REST service recently added"
$ws.Range("H31").Value = "https://droughtmonitor.unl.edu/Data/DataDownload/WeeksInDrought.aspx"

# Hyperlink for the info link cell
$ws.Hyperlinks.Add($ws.Range("H31"), "https://droughtmonitor.unl.edu/Data/DataDownload/WeeksInDrought.aspx") | Out-Null

# Re-apply the built-in Hyperlink style so H31 matches the other Info Link cells
$ws.Range("H31").Style = "Hyperlink"

# --- Sheet view state --------------------------------------------------------
# Move the active selection / scroll position as recorded in the workbook view
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("G26").Select() | Out-Null
